$wb = $excel.ActiveWorkbook

# --- "readme" sheet: reorder the JobNo / sheet_name / Author columns ------
# Old layout: A=index, B=JobNo, C=sheet_name, D=Author, E=Date
# New layout: A=index, B=Author, C=JobNo,     D=sheet_name, E=Date
$ws = $wb.Worksheets.Item("readme")

# Update the table header row (this also rewrites the ListObject's
# tableColumn names automatically).
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "JobNo"
$ws.Range("D1").Value = "sheet_name"

# Sheet names in row order (rows 2..12), used for the new "sheet_name"
# column (D).
$sheetNames = @(
    "Project Information",
    "Criterion Definitions",
    "Results, Air Speed 0.1",
    "Results, Air Speed 0.15",
    "Results, Air Speed 0.2",
    "Results, Air Speed 0.3",
    "Results, Air Speed 0.4",
    "Results, Air Speed 0.5",
    "Results, Air Speed 0.6",
    "Results, Air Speed 0.7",
    "Results, Air Speed 0.8"
)

for ($i = 0; $i -lt $sheetNames.Length; $i++) {
    $row = 2 + $i
    $ws.Range("B$row").Value = "jovyan"
    $ws.Range("C$row").Value = "/c/e"
    $ws.Range("D$row").Value = $sheetNames[$i]
}

# --- "Project Information" sheet: refresh the analysis timestamp ---------
$ws2 = $wb.Worksheets.Item("Project Information")
$ws2.Range("B12").Value = "2022-06-15 15:57:18.086952"
